# Alpha3F-HW25 -- updated notebook, reran simulation
#
# Two brand-new sample rows ("Holden" and "Rizzie Spiral") land right
# after "Spiral5"; the "Thomas Hex" sample is renamed "Matthies Hex";
# and the whole averaging simulation is rerun, so every HKL row's
# computed ratios (columns C:T) refresh with freshly simulated
# numbers. The table grows from 29 to 31 data rows (A1:T29 -> A1:T31).
#
# The upstream tool (a Jupyter notebook) regenerates this sheet from
# scratch on every run, so rather than shuffling existing rows in
# place, each row is rewritten directly with its final value -- the
# same end state Excel would show after a "Save As" of the rerun.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4: Holden (new)
$ws.Cells.Item(4,1).Value = 2
$ws.Cells.Item(4,2).Value = 'Holden'
$ws.Cells.Item(4,3).Value = 1.01209127276858
$ws.Cells.Item(4,4).Value = 0.9731341077459776
$ws.Cells.Item(4,5).Value = 1.042368299902967
$ws.Cells.Item(4,6).Value = 1.002112216537267
$ws.Cells.Item(4,7).Value = 0.9539394133812794
$ws.Cells.Item(4,8).Value = 1.002112216537267
$ws.Cells.Item(4,9).Value = 0.9539394133812794
$ws.Cells.Item(4,10).Value = 1.016197877769689
$ws.Cells.Item(4,11).Value = 1.002112216537267
$ws.Cells.Item(4,12).Value = 1.016197877769689
$ws.Cells.Item(4,13).Value = 0.985068645575484
$ws.Cells.Item(4,14).Value = 0.985068645575484
$ws.Cells.Item(4,15).Value = 0.981090466298982
$ws.Cells.Item(4,16).Value = 0.9907498358960782
$ws.Cells.Item(4,17).Value = 0.9907498358960782
$ws.Cells.Item(4,18).Value = 0.9935904310563753
$ws.Cells.Item(4,19).Value = 0.9935904310563753
$ws.Cells.Item(4,20).Value = 0.9999738646842932

# Row 5: Rizzie Spiral (new)
$ws.Cells.Item(5,1).Value = 3
$ws.Cells.Item(5,2).Value = 'Rizzie Spiral'
$ws.Cells.Item(5,3).Value = 1.094311840689
$ws.Cells.Item(5,4).Value = 0.8015858300408976
$ws.Cells.Item(5,5).Value = 1.12969791767187
$ws.Cells.Item(5,6).Value = 1.18953405722407
$ws.Cells.Item(5,7).Value = 0.5942908985687556
$ws.Cells.Item(5,8).Value = 1.18953405722407
$ws.Cells.Item(5,9).Value = 0.5942908985687556
$ws.Cells.Item(5,10).Value = 1.055125744823871
$ws.Cells.Item(5,11).Value = 1.18953405722407
$ws.Cells.Item(5,12).Value = 1.055125744823871
$ws.Cells.Item(5,13).Value = 0.8247083216963135
$ws.Cells.Item(5,14).Value = 0.8247083216963135
$ws.Cells.Item(5,15).Value = 0.8170008244778414
$ws.Cells.Item(5,16).Value = 0.9463169002055659
$ws.Cells.Item(5,17).Value = 0.9463169002055659
$ws.Cells.Item(5,18).Value = 1.007121189460192
$ws.Cells.Item(5,19).Value = 1.007121189460192
$ws.Cells.Item(5,20).Value = 0.9774243815030775

# Row 6: RotRing OmegaMax-90
$ws.Cells.Item(6,1).Value = 4
$ws.Cells.Item(6,2).Value = 'RotRing OmegaMax-90'
$ws.Cells.Item(6,3).Value = 1.029860053243006
$ws.Cells.Item(6,4).Value = 0.9306570305502204
$ws.Cells.Item(6,5).Value = 1.115163487510936
$ws.Cells.Item(6,6).Value = 1.005755844614242
$ws.Cells.Item(6,7).Value = 0.8902581062539257
$ws.Cells.Item(6,8).Value = 1.005755844614242
$ws.Cells.Item(6,9).Value = 0.8902581062539257
$ws.Cells.Item(6,10).Value = 1.039779476010362
$ws.Cells.Item(6,11).Value = 1.005755844614242
$ws.Cells.Item(6,12).Value = 1.039779476010362
$ws.Cells.Item(6,13).Value = 0.9650187911321436
$ws.Cells.Item(6,14).Value = 0.9650187911321436
$ws.Cells.Item(6,15).Value = 0.9535648709381692
$ws.Cells.Item(6,16).Value = 0.9785978089595098
$ws.Cells.Item(6,17).Value = 0.9785978089595098
$ws.Cells.Item(6,18).Value = 0.9853873178731928
$ws.Cells.Item(6,19).Value = 0.9853873178731928
$ws.Cells.Item(6,20).Value = 1.001912333030449

# Row 7: Equal Angle
$ws.Cells.Item(7,1).Value = 5
$ws.Cells.Item(7,2).Value = 'Equal Angle'
$ws.Cells.Item(7,3).Value = 1.01583789538905
$ws.Cells.Item(7,4).Value = 0.9652271944092247
$ws.Cells.Item(7,5).Value = 1.050289408804033
$ws.Cells.Item(7,6).Value = 1.006764224257925
$ws.Cells.Item(7,7).Value = 0.9383261708573483
$ws.Cells.Item(7,8).Value = 1.006764224257925
$ws.Cells.Item(7,9).Value = 0.9383261708573483
$ws.Cells.Item(7,10).Value = 1.019571911945244
$ws.Cells.Item(7,11).Value = 1.006764224257925
$ws.Cells.Item(7,12).Value = 1.019571911945244
$ws.Cells.Item(7,13).Value = 0.9789490414012962
$ws.Cells.Item(7,14).Value = 0.9789490414012962
$ws.Cells.Item(7,15).Value = 0.974375092403939
$ws.Cells.Item(7,16).Value = 0.9882207690201725
$ws.Cells.Item(7,17).Value = 0.9882207690201725
$ws.Cells.Item(7,18).Value = 0.9928566328296107
$ws.Cells.Item(7,19).Value = 0.9928566328296107
$ws.Cells.Item(7,20).Value = 0.9993361342771375

# Row 8: Tilt Rotate
$ws.Cells.Item(8,1).Value = 6
$ws.Cells.Item(8,2).Value = 'Tilt Rotate'
$ws.Cells.Item(8,3).Value = 1.04926171772109
$ws.Cells.Item(8,4).Value = 0.8916512474684229
$ws.Cells.Item(8,5).Value = 1.178462022008578
$ws.Cells.Item(8,6).Value = 0.9975603262270364
$ws.Cells.Item(8,7).Value = 0.8134353594606664
$ws.Cells.Item(8,8).Value = 0.9975603262270364
$ws.Cells.Item(8,9).Value = 0.8134353594606664
$ws.Cells.Item(8,10).Value = 1.070537972377412
$ws.Cells.Item(8,11).Value = 0.9975603262270364
$ws.Cells.Item(8,12).Value = 1.070537972377412
$ws.Cells.Item(8,13).Value = 0.941986665919039
$ws.Cells.Item(8,14).Value = 0.941986665919039
$ws.Cells.Item(8,15).Value = 0.9252081931021671
$ws.Cells.Item(8,16).Value = 0.9605112193550381
$ws.Cells.Item(8,17).Value = 0.9605112193550381
$ws.Cells.Item(8,18).Value = 0.9697734960730378
$ws.Cells.Item(8,19).Value = 0.9697734960730378
$ws.Cells.Item(8,20).Value = 1.000151440877201

# Row 9: CLR
$ws.Cells.Item(9,1).Value = 7
$ws.Cells.Item(9,2).Value = 'CLR'
$ws.Cells.Item(9,3).Value = 1.002023901731235
$ws.Cells.Item(9,4).Value = 0.9958271379633181
$ws.Cells.Item(9,5).Value = 1.00294605781747
$ws.Cells.Item(9,6).Value = 1.003615148924999
$ws.Cells.Item(9,7).Value = 0.9911263220042859
$ws.Cells.Item(9,8).Value = 1.003615148924999
$ws.Cells.Item(9,9).Value = 0.9911263220042859
$ws.Cells.Item(9,10).Value = 1.001369065859894
$ws.Cells.Item(9,11).Value = 1.003615148924999
$ws.Cells.Item(9,12).Value = 1.001369065859894
$ws.Cells.Item(9,13).Value = 0.9962476939320899
$ws.Cells.Item(9,14).Value = 0.9962476939320899
$ws.Cells.Item(9,15).Value = 0.9961075086091661
$ws.Cells.Item(9,16).Value = 0.9987035122630598
$ws.Cells.Item(9,17).Value = 0.9987035122630598
$ws.Cells.Item(9,18).Value = 0.9999314214285446
$ws.Cells.Item(9,19).Value = 0.9999314214285446
$ws.Cells.Item(9,20).Value = 0.999484605716867

# Row 10: Rizzie Hex
$ws.Cells.Item(10,1).Value = 8
$ws.Cells.Item(10,2).Value = 'Rizzie Hex'
$ws.Cells.Item(10,3).Value = 1.000175584650901
$ws.Cells.Item(10,4).Value = 0.9996953136296278
$ws.Cells.Item(10,5).Value = 0.9994981217320993
$ws.Cells.Item(10,6).Value = 1.000904345709062
$ws.Cells.Item(10,7).Value = 0.9990476075662951
$ws.Cells.Item(10,8).Value = 1.000904345709062
$ws.Cells.Item(10,9).Value = 0.9990476075662951
$ws.Cells.Item(10,10).Value = 0.999875682363654
$ws.Cells.Item(10,11).Value = 1.000904345709062
$ws.Cells.Item(10,12).Value = 0.999875682363654
$ws.Cells.Item(10,13).Value = 0.9994616449649745
$ws.Cells.Item(10,14).Value = 0.9994616449649745
$ws.Cells.Item(10,15).Value = 0.999539534519859
$ws.Cells.Item(10,16).Value = 0.9999425452130035
$ws.Cells.Item(10,17).Value = 0.9999425452130036
$ws.Cells.Item(10,18).Value = 1.000182995337018
$ws.Cells.Item(10,19).Value = 1.000182995337018
$ws.Cells.Item(10,20).Value = 0.9998661092752732

# Row 11: Matthies Hex (renamed from Thomas Hex)
$ws.Cells.Item(11,1).Value = 9
$ws.Cells.Item(11,2).Value = 'Matthies Hex'
$ws.Cells.Item(11,3).Value = 1.002967167807092
$ws.Cells.Item(11,4).Value = 0.9939135626727936
$ws.Cells.Item(11,5).Value = 1.004957789339782
$ws.Cells.Item(11,6).Value = 1.004514637560676
$ws.Cells.Item(11,7).Value = 0.9870244538872426
$ws.Cells.Item(11,8).Value = 1.004514637560676
$ws.Cells.Item(11,9).Value = 0.9870244538872426
$ws.Cells.Item(11,10).Value = 1.002330348206932
$ws.Cells.Item(11,11).Value = 1.004514637560676
$ws.Cells.Item(11,12).Value = 1.002330348206932
$ws.Cells.Item(11,13).Value = 0.9946774010470874
$ws.Cells.Item(11,14).Value = 0.9946774010470874
$ws.Cells.Item(11,15).Value = 0.9944227882556561
$ws.Cells.Item(11,16).Value = 0.9979564798849502
$ws.Cells.Item(11,17).Value = 0.9979564798849502
$ws.Cells.Item(11,18).Value = 0.9995960193038815
$ws.Cells.Item(11,19).Value = 0.9995960193038815
$ws.Cells.Item(11,20).Value = 0.9992846599124197

# Row 12: Tilt Rotate_Partial
$ws.Cells.Item(12,1).Value = 10
$ws.Cells.Item(12,2).Value = 'Tilt Rotate_Partial'
$ws.Cells.Item(12,3).Value = 1.050022242112039
$ws.Cells.Item(12,4).Value = 0.8901285600602519
$ws.Cells.Item(12,5).Value = 1.181723565125295
$ws.Cells.Item(12,6).Value = 0.9964183924769775
$ws.Cells.Item(12,7).Value = 0.8104593637871472
$ws.Cells.Item(12,8).Value = 0.9964183924769775
$ws.Cells.Item(12,9).Value = 0.8104593637871472
$ws.Cells.Item(12,10).Value = 1.072081400736139
$ws.Cells.Item(12,11).Value = 0.9964183924769775
$ws.Cells.Item(12,12).Value = 1.072081400736139
$ws.Cells.Item(12,13).Value = 0.9412703822616433
$ws.Cells.Item(12,14).Value = 0.9412703822616433
$ws.Cells.Item(12,15).Value = 0.9242231081945128
$ws.Cells.Item(12,16).Value = 0.9596530523334214
$ws.Cells.Item(12,17).Value = 0.9596530523334215
$ws.Cells.Item(12,18).Value = 0.9688443873693104
$ws.Cells.Item(12,19).Value = 0.9688443873693104
$ws.Cells.Item(12,20).Value = 1.000138920716308

# Row 13: RotRing OmegaMax-60
$ws.Cells.Item(13,1).Value = 11
$ws.Cells.Item(13,2).Value = 'RotRing OmegaMax-60'
$ws.Cells.Item(13,3).Value = 1.030145054151836
$ws.Cells.Item(13,4).Value = 0.9305910766820134
$ws.Cells.Item(13,5).Value = 1.11129016368297
$ws.Cells.Item(13,6).Value = 1.008907961864615
$ws.Cells.Item(13,7).Value = 0.8876539919539809
$ws.Cells.Item(13,8).Value = 1.008907961864615
$ws.Cells.Item(13,9).Value = 0.8876539919539809
$ws.Cells.Item(13,10).Value = 1.038884594350739
$ws.Cells.Item(13,11).Value = 1.008907961864615
$ws.Cells.Item(13,12).Value = 1.038884594350739
$ws.Cells.Item(13,13).Value = 0.9632692931523597
$ws.Cells.Item(13,14).Value = 0.9632692931523597
$ws.Cells.Item(13,15).Value = 0.9523765543289109
$ws.Cells.Item(13,16).Value = 0.9784821827231115
$ws.Cells.Item(13,17).Value = 0.9784821827231115
$ws.Cells.Item(13,18).Value = 0.9860886275084875
$ws.Cells.Item(13,19).Value = 0.9860886275084875
$ws.Cells.Item(13,20).Value = 1.001245473781026

# Row 14: Equal Angle_Partial
$ws.Cells.Item(14,1).Value = 12
$ws.Cells.Item(14,2).Value = 'Equal Angle_Partial'
$ws.Cells.Item(14,3).Value = 1.016015130168423
$ws.Cells.Item(14,4).Value = 0.9648286567157899
$ws.Cells.Item(14,5).Value = 1.056334868378946
$ws.Cells.Item(14,6).Value = 1.000955043084212
$ws.Cells.Item(14,7).Value = 0.9386744326736818
$ws.Cells.Item(14,8).Value = 1.000955043084212
$ws.Cells.Item(14,9).Value = 0.9386744326736818
$ws.Cells.Item(14,10).Value = 1.022212691231578
$ws.Cells.Item(14,11).Value = 1.000955043084212
$ws.Cells.Item(14,12).Value = 1.022212691231578
$ws.Cells.Item(14,13).Value = 0.9804435619526299
$ws.Cells.Item(14,14).Value = 0.9804435619526299
$ws.Cells.Item(14,15).Value = 0.9752385935403499
$ws.Cells.Item(14,16).Value = 0.987280722329824
$ws.Cells.Item(14,17).Value = 0.987280722329824
$ws.Cells.Item(14,18).Value = 0.990699302518421
$ws.Cells.Item(14,19).Value = 0.990699302518421
$ws.Cells.Item(14,20).Value = 0.9998368037087718

# Row 15: Rizzie Hex_Partial
$ws.Cells.Item(15,1).Value = 13
$ws.Cells.Item(15,2).Value = 'Rizzie Hex_Partial'
$ws.Cells.Item(15,3).Value = 0.9901318950068353
$ws.Cells.Item(15,4).Value = 1.022213070888607
$ws.Cells.Item(15,5).Value = 0.9721911552389385
$ws.Cells.Item(15,6).Value = 0.989899920014533
$ws.Cells.Item(15,7).Value = 1.038497272460191
$ws.Cells.Item(15,8).Value = 0.989899920014533
$ws.Cells.Item(15,9).Value = 1.038497272460191
$ws.Cells.Item(15,10).Value = 0.9902273579528168
$ws.Cells.Item(15,11).Value = 0.989899920014533
$ws.Cells.Item(15,12).Value = 0.9902273579528168
$ws.Cells.Item(15,13).Value = 1.014362315206504
$ws.Cells.Item(15,14).Value = 1.014362315206504
$ws.Cells.Item(15,15).Value = 1.016979233767205
$ws.Cells.Item(15,16).Value = 1.006208183475847
$ws.Cells.Item(15,17).Value = 1.006208183475847
$ws.Cells.Item(15,18).Value = 1.002131117610518
$ws.Cells.Item(15,19).Value = 1.002131117610518
$ws.Cells.Item(15,20).Value = 1.000526778593654

# Row 16: ND Single
$ws.Cells.Item(16,1).Value = 14
$ws.Cells.Item(16,2).Value = 'ND Single'
$ws.Cells.Item(16,3).Value = 1.0863871
$ws.Cells.Item(16,4).Value = 0.8102074900000015
$ws.Cells.Item(16,5).Value = 1.315735799999999
$ws.Cells.Item(16,6).Value = 0.9917943100000015
$ws.Cells.Item(16,7).Value = 0.6734045099999993
$ws.Cells.Item(16,8).Value = 0.9917943100000015
$ws.Cells.Item(16,9).Value = 0.6734045099999993
$ws.Cells.Item(16,10).Value = 1.125314099999998
$ws.Cells.Item(16,11).Value = 0.9917943100000015
$ws.Cells.Item(16,12).Value = 1.125314099999998
$ws.Cells.Item(16,13).Value = 0.8993593049999988
$ws.Cells.Item(16,14).Value = 0.8993593049999988
$ws.Cells.Item(16,15).Value = 0.8696420333333331
$ws.Cells.Item(16,16).Value = 0.9301709733333331
$ws.Cells.Item(16,17).Value = 0.9301709733333331
$ws.Cells.Item(16,18).Value = 0.9455768075000002
$ws.Cells.Item(16,19).Value = 0.9455768075000002
$ws.Cells.Item(16,20).Value = 1.000473885

# Row 17: RD Single
$ws.Cells.Item(17,1).Value = 15
$ws.Cells.Item(17,2).Value = 'RD Single'
$ws.Cells.Item(17,3).Value = 1.1014624
$ws.Cells.Item(17,4).Value = 0.7938064500000001
$ws.Cells.Item(17,5).Value = 0.9618347200000001
$ws.Cells.Item(17,6).Value = 1.3679559
$ws.Cells.Item(17,7).Value = 0.5229061699999999
$ws.Cells.Item(17,8).Value = 1.3679559
$ws.Cells.Item(17,9).Value = 0.5229061699999999
$ws.Cells.Item(17,10).Value = 0.99179431
$ws.Cells.Item(17,11).Value = 1.3679559
$ws.Cells.Item(17,12).Value = 0.99179431
$ws.Cells.Item(17,13).Value = 0.75735024
$ws.Cells.Item(17,14).Value = 0.75735024
$ws.Cells.Item(17,15).Value = 0.76950231
$ws.Cells.Item(17,16).Value = 0.9608854599999997
$ws.Cells.Item(17,17).Value = 0.9608854599999997
$ws.Cells.Item(17,18).Value = 1.06265307
$ws.Cells.Item(17,19).Value = 1.06265307
$ws.Cells.Item(17,20).Value = 0.9566266583333332

# Row 18: TD Single
$ws.Cells.Item(18,1).Value = 16
$ws.Cells.Item(18,2).Value = 'TD Single'
$ws.Cells.Item(18,3).Value = 1.2125577
$ws.Cells.Item(18,4).Value = 0.49189623
$ws.Cells.Item(18,5).Value = 1.9473815
$ws.Cells.Item(18,6).Value = 0.9618347200000001
$ws.Cells.Item(18,7).Value = 0.25020352
$ws.Cells.Item(18,8).Value = 0.9618347200000001
$ws.Cells.Item(18,9).Value = 0.25020352
$ws.Cells.Item(18,10).Value = 1.3157358
$ws.Cells.Item(18,11).Value = 0.9618347200000001
$ws.Cells.Item(18,12).Value = 1.3157358
$ws.Cells.Item(18,13).Value = 0.78296966
$ws.Cells.Item(18,14).Value = 0.78296966
$ws.Cells.Item(18,15).Value = 0.6859451833333333
$ws.Cells.Item(18,16).Value = 0.8425913466666667
$ws.Cells.Item(18,17).Value = 0.8425913466666666
$ws.Cells.Item(18,18).Value = 0.87240219
$ws.Cells.Item(18,19).Value = 0.87240219
$ws.Cells.Item(18,20).Value = 1.029934911666667

# Row 19: Morris Single
$ws.Cells.Item(19,1).Value = 17
$ws.Cells.Item(19,2).Value = 'Morris Single'
$ws.Cells.Item(19,3).Value = 1.012886
$ws.Cells.Item(19,4).Value = 0.98424247
$ws.Cells.Item(19,5).Value = 0.90037857
$ws.Cells.Item(19,6).Value = 1.1082556
$ws.Cells.Item(19,7).Value = 0.91308462
$ws.Cells.Item(19,8).Value = 1.1082556
$ws.Cells.Item(19,9).Value = 0.91308462
$ws.Cells.Item(19,10).Value = 0.97363919
$ws.Cells.Item(19,11).Value = 1.1082556
$ws.Cells.Item(19,12).Value = 0.97363919
$ws.Cells.Item(19,13).Value = 0.9433619049999999
$ws.Cells.Item(19,14).Value = 0.9433619049999999
$ws.Cells.Item(19,15).Value = 0.9569887599999999
$ws.Cells.Item(19,16).Value = 0.99832647
$ws.Cells.Item(19,17).Value = 0.99832647
$ws.Cells.Item(19,18).Value = 1.0258087525
$ws.Cells.Item(19,19).Value = 1.0258087525
$ws.Cells.Item(19,20).Value = 0.982081075

# Row 20: Ring Perpendicular to ND
$ws.Cells.Item(20,1).Value = 18
$ws.Cells.Item(20,2).Value = 'Ring Perpendicular to ND'
$ws.Cells.Item(20,3).Value = 1.036389898630137
$ws.Cells.Item(20,4).Value = 0.9253719497260274
$ws.Cells.Item(20,5).Value = 1.057229300821918
$ws.Cells.Item(20,6).Value = 1.059140436986302
$ws.Cells.Item(20,7).Value = 0.8402270360273973
$ws.Cells.Item(20,8).Value = 1.059140436986302
$ws.Cells.Item(20,9).Value = 0.8402270360273973
$ws.Cells.Item(20,10).Value = 1.027027528904109
$ws.Cells.Item(20,11).Value = 1.059140436986302
$ws.Cells.Item(20,12).Value = 1.027027528904109
$ws.Cells.Item(20,13).Value = 0.9336272824657533
$ws.Cells.Item(20,14).Value = 0.9336272824657533
$ws.Cells.Item(20,15).Value = 0.9308755048858447
$ws.Cells.Item(20,16).Value = 0.9754650006392694
$ws.Cells.Item(20,17).Value = 0.9754650006392694
$ws.Cells.Item(20,18).Value = 0.9963838597260275
$ws.Cells.Item(20,19).Value = 0.9963838597260275
$ws.Cells.Item(20,20).Value = 0.9908976918493151

# Row 21: Ring Perpendicular to RD
$ws.Cells.Item(21,1).Value = 19
$ws.Cells.Item(21,2).Value = 'Ring Perpendicular to RD'
$ws.Cells.Item(21,3).Value = 1.032922055789474
$ws.Cells.Item(21,4).Value = 0.9144072342105264
$ws.Cells.Item(21,5).Value = 1.203647928421053
$ws.Cells.Item(21,6).Value = 0.9584451015789475
$ws.Cells.Item(21,7).Value = 0.9028523952631577
$ws.Cells.Item(21,8).Value = 0.9584451015789475
$ws.Cells.Item(21,9).Value = 0.9028523952631577
$ws.Cells.Item(21,10).Value = 1.063570988421053
$ws.Cells.Item(21,11).Value = 0.9584451015789475
$ws.Cells.Item(21,12).Value = 1.063570988421053
$ws.Cells.Item(21,13).Value = 0.9832116918421052
$ws.Cells.Item(21,14).Value = 0.9832116918421052
$ws.Cells.Item(21,15).Value = 0.9602768726315789
$ws.Cells.Item(21,16).Value = 0.9749561617543859
$ws.Cells.Item(21,17).Value = 0.9749561617543859
$ws.Cells.Item(21,18).Value = 0.9708283967105262
$ws.Cells.Item(21,19).Value = 0.9708283967105262
$ws.Cells.Item(21,20).Value = 1.012640950614035

# Row 22: Ring Perpendicular to TD
$ws.Cells.Item(22,1).Value = 20
$ws.Cells.Item(22,2).Value = 'Ring Perpendicular to TD'
$ws.Cells.Item(22,3).Value = 1.093924742105263
$ws.Cells.Item(22,4).Value = 0.8020069699999999
$ws.Cells.Item(22,5).Value = 1.13878526
$ws.Cells.Item(22,6).Value = 1.179875116842106
$ws.Cells.Item(22,7).Value = 0.5981553399999998
$ws.Cells.Item(22,8).Value = 1.179875116842106
$ws.Cells.Item(22,9).Value = 0.5981553399999998
$ws.Cells.Item(22,10).Value = 1.058554216842105
$ws.Cells.Item(22,11).Value = 1.179875116842106
$ws.Cells.Item(22,12).Value = 1.058554216842105
$ws.Cells.Item(22,13).Value = 0.8283547784210525
$ws.Cells.Item(22,14).Value = 0.8283547784210525
$ws.Cells.Item(22,15).Value = 0.8195721756140349
$ws.Cells.Item(22,16).Value = 0.9455282245614036
$ws.Cells.Item(22,17).Value = 0.9455282245614036
$ws.Cells.Item(22,18).Value = 1.004114947631579
$ws.Cells.Item(22,19).Value = 1.004114947631579
$ws.Cells.Item(22,20).Value = 0.9785502742982456

# Row 23: OffsetFTD
$ws.Cells.Item(23,1).Value = 21
$ws.Cells.Item(23,2).Value = 'OffsetFTD'
$ws.Cells.Item(23,3).Value = 0.9237120272003099
$ws.Cells.Item(23,4).Value = 1.155687464767944
$ws.Cells.Item(23,5).Value = 0.8891843764791993
$ws.Cells.Item(23,6).Value = 0.8700956238553862
$ws.Cells.Item(23,7).Value = 1.335061443235451
$ws.Cells.Item(23,8).Value = 0.8700956238553862
$ws.Cells.Item(23,9).Value = 1.335061443235451
$ws.Cells.Item(23,10).Value = 0.9457763916802239
$ws.Cells.Item(23,11).Value = 0.8700956238553862
$ws.Cells.Item(23,12).Value = 0.9457763916802239
$ws.Cells.Item(23,13).Value = 1.140418917457837
$ws.Cells.Item(23,14).Value = 1.140418917457837
$ws.Cells.Item(23,15).Value = 1.145508433227873
$ws.Cells.Item(23,16).Value = 1.050311152923687
$ws.Cells.Item(23,17).Value = 1.050311152923687
$ws.Cells.Item(23,18).Value = 1.005257270656612
$ws.Cells.Item(23,19).Value = 1.005257270656612
$ws.Cells.Item(23,20).Value = 1.019919554536419

# Row 24: OffsetATD
$ws.Cells.Item(24,1).Value = 22
$ws.Cells.Item(24,2).Value = 'OffsetATD'
$ws.Cells.Item(24,3).Value = 0.9847977284212999
$ws.Cells.Item(24,4).Value = 1.04199403505895
$ws.Cells.Item(24,5).Value = 0.9161628456874087
$ws.Cells.Item(24,6).Value = 0.9994301372983533
$ws.Cells.Item(24,7).Value = 1.042647621188965
$ws.Cells.Item(24,8).Value = 0.9994301372983533
$ws.Cells.Item(24,9).Value = 1.042647621188965
$ws.Cells.Item(24,10).Value = 0.9787761668046858
$ws.Cells.Item(24,11).Value = 0.9994301372983533
$ws.Cells.Item(24,12).Value = 0.9787761668046858
$ws.Cells.Item(24,13).Value = 1.010711893996826
$ws.Cells.Item(24,14).Value = 1.010711893996826
$ws.Cells.Item(24,15).Value = 1.021139274350867
$ws.Cells.Item(24,16).Value = 1.006951308430668
$ws.Cells.Item(24,17).Value = 1.006951308430668
$ws.Cells.Item(24,18).Value = 1.005071015647589
$ws.Cells.Item(24,19).Value = 1.005071015647589
$ws.Cells.Item(24,20).Value = 0.9939680890766104

# Row 25: OffsetF45
$ws.Cells.Item(25,1).Value = 23
$ws.Cells.Item(25,2).Value = 'OffsetF45'
$ws.Cells.Item(25,3).Value = 0.9789606011489488
$ws.Cells.Item(25,4).Value = 1.05017208908841
$ws.Cells.Item(25,5).Value = 0.9327484584809496
$ws.Cells.Item(25,6).Value = 0.9762328983501093
$ws.Cells.Item(25,7).Value = 1.077877451537533
$ws.Cells.Item(25,8).Value = 0.9762328983501093
$ws.Cells.Item(25,9).Value = 1.077877451537533
$ws.Cells.Item(25,10).Value = 0.9800831147905154
$ws.Cells.Item(25,11).Value = 0.9762328983501093
$ws.Cells.Item(25,12).Value = 0.9800831147905154
$ws.Cells.Item(25,13).Value = 1.028980283164024
$ws.Cells.Item(25,14).Value = 1.028980283164024
$ws.Cells.Item(25,15).Value = 1.036044218472153
$ws.Cells.Item(25,16).Value = 1.011397821559386
$ws.Cells.Item(25,17).Value = 1.011397821559386
$ws.Cells.Item(25,18).Value = 1.002606590757067
$ws.Cells.Item(25,19).Value = 1.002606590757067
$ws.Cells.Item(25,20).Value = 0.9993457688994112

# Row 26: OffsetA45
$ws.Cells.Item(26,1).Value = 24
$ws.Cells.Item(26,2).Value = 'OffsetA45'
$ws.Cells.Item(26,3).Value = 0.9926053784202922
$ws.Cells.Item(26,4).Value = 1.011857164849821
$ws.Cells.Item(26,5).Value = 0.9900919982056045
$ws.Cells.Item(26,6).Value = 0.9988460005791482
$ws.Cells.Item(26,7).Value = 1.036008032146925
$ws.Cells.Item(26,8).Value = 0.9988460005791482
$ws.Cells.Item(26,9).Value = 1.036008032146925
$ws.Cells.Item(26,10).Value = 0.9900372152319273
$ws.Cells.Item(26,11).Value = 0.9988460005791482
$ws.Cells.Item(26,12).Value = 0.9900372152319273
$ws.Cells.Item(26,13).Value = 1.013022623689426
$ws.Cells.Item(26,14).Value = 1.013022623689426
$ws.Cells.Item(26,15).Value = 1.012634137409558
$ws.Cells.Item(26,16).Value = 1.008297082652667
$ws.Cells.Item(26,17).Value = 1.008297082652667
$ws.Cells.Item(26,18).Value = 1.005934312134287
$ws.Cells.Item(26,19).Value = 1.005934312134287
$ws.Cells.Item(26,20).Value = 1.00324096490562

# Row 27: OffsetFRD
$ws.Cells.Item(27,1).Value = 25
$ws.Cells.Item(27,2).Value = 'OffsetFRD'
$ws.Cells.Item(27,3).Value = 1.076831202290409
$ws.Cells.Item(27,4).Value = 0.8426804539050847
$ws.Cells.Item(27,5).Value = 1.040066695130711
$ws.Cells.Item(27,6).Value = 1.209199205433904
$ws.Cells.Item(27,7).Value = 0.6516794280098702
$ws.Cells.Item(27,8).Value = 1.209199205433904
$ws.Cells.Item(27,9).Value = 0.6516794280098702
$ws.Cells.Item(27,10).Value = 1.022358760569408
$ws.Cells.Item(27,11).Value = 1.209199205433904
$ws.Cells.Item(27,12).Value = 1.022358760569408
$ws.Cells.Item(27,13).Value = 0.8370190942896394
$ws.Cells.Item(27,14).Value = 0.8370190942896394
$ws.Cells.Item(27,15).Value = 0.8389062141614545
$ws.Cells.Item(27,16).Value = 0.9610791313377275
$ws.Cells.Item(27,17).Value = 0.9610791313377275
$ws.Cells.Item(27,18).Value = 1.023109149861771
$ws.Cells.Item(27,19).Value = 1.023109149861771
$ws.Cells.Item(27,20).Value = 0.9738026242232311

# Row 28: OffsetARD
$ws.Cells.Item(28,1).Value = 26
$ws.Cells.Item(28,2).Value = 'OffsetARD'
$ws.Cells.Item(28,3).Value = 1.01624753496696
$ws.Cells.Item(28,4).Value = 0.9624088333558491
$ws.Cells.Item(28,5).Value = 1.106335847993797
$ws.Cells.Item(28,6).Value = 0.9557657070441085
$ws.Cells.Item(28,7).Value = 0.9475848549883495
$ws.Cells.Item(28,8).Value = 0.9557657070441085
$ws.Cells.Item(28,9).Value = 0.9475848549883495
$ws.Cells.Item(28,10).Value = 1.041137176235208
$ws.Cells.Item(28,11).Value = 0.9557657070441085
$ws.Cells.Item(28,12).Value = 1.041137176235208
$ws.Cells.Item(28,13).Value = 0.9943610156117786
$ws.Cells.Item(28,14).Value = 0.9943610156117786
$ws.Cells.Item(28,15).Value = 0.9837102881931354
$ws.Cells.Item(28,16).Value = 0.9814959127558885
$ws.Cells.Item(28,17).Value = 0.9814959127558884
$ws.Cells.Item(28,18).Value = 0.9750633613279435
$ws.Cells.Item(28,19).Value = 0.9750633613279435
$ws.Cells.Item(28,20).Value = 1.004913325764045

# Row 29: Gaussian Quadrature
$ws.Cells.Item(29,1).Value = 27
$ws.Cells.Item(29,2).Value = 'Gaussian Quadrature'
$ws.Cells.Item(29,3).Value = 1.00463166672915
$ws.Cells.Item(29,4).Value = 0.9871814698368729
$ws.Cells.Item(29,5).Value = 1.013743867571337
$ws.Cells.Item(29,6).Value = 1.012021657389473
$ws.Cells.Item(29,7).Value = 0.9870565266743779
$ws.Cells.Item(29,8).Value = 1.012021657389473
$ws.Cells.Item(29,9).Value = 0.9870565266743779
$ws.Cells.Item(29,10).Value = 1.001590508655779
$ws.Cells.Item(29,11).Value = 1.012021657389473
$ws.Cells.Item(29,12).Value = 1.001590508655779
$ws.Cells.Item(29,13).Value = 0.9943235176650786
$ws.Cells.Item(29,14).Value = 0.9943235176650786
$ws.Cells.Item(29,15).Value = 0.9919428350556768
$ws.Cells.Item(29,16).Value = 1.00022289757321
$ws.Cells.Item(29,17).Value = 1.00022289757321
$ws.Cells.Item(29,18).Value = 1.003172587527276
$ws.Cells.Item(29,19).Value = 1.003172587527276
$ws.Cells.Item(29,20).Value = 1.001037616142832

# Row 30: Michael-CCHex (new row)
$ws.Cells.Item(30,1).Value = 28
$ws.Cells.Item(30,2).Value = 'Michael-CCHex'
$ws.Cells.Item(30,3).Value = 1.009844347310388
$ws.Cells.Item(30,4).Value = 0.9797353545962898
$ws.Cells.Item(30,5).Value = 1.013527238219397
$ws.Cells.Item(30,6).Value = 1.018391404170891
$ws.Cells.Item(30,7).Value = 0.9565085723005406
$ws.Cells.Item(30,8).Value = 1.018391404170891
$ws.Cells.Item(30,9).Value = 0.9565085723005406
$ws.Cells.Item(30,10).Value = 1.006327043632693
$ws.Cells.Item(30,11).Value = 1.018391404170891
$ws.Cells.Item(30,12).Value = 1.006327043632693
$ws.Cells.Item(30,13).Value = 0.9814178079666169
$ws.Cells.Item(30,14).Value = 0.9814178079666169
$ws.Cells.Item(30,15).Value = 0.980856990176508
$ws.Cells.Item(30,16).Value = 0.9937423400347084
$ws.Cells.Item(30,17).Value = 0.9937423400347084
$ws.Cells.Item(30,18).Value = 0.999904606068754
$ws.Cells.Item(30,19).Value = 0.999904606068754
$ws.Cells.Item(30,20).Value = 0.9973889933717

# Row 31: Michael-SNHex (new row)
$ws.Cells.Item(31,1).Value = 29
$ws.Cells.Item(31,2).Value = 'Michael-SNHex'
$ws.Cells.Item(31,3).Value = 0.9934706138223061
$ws.Cells.Item(31,4).Value = 1.014605155172839
$ws.Cells.Item(31,5).Value = 0.9417775132367981
$ws.Cells.Item(31,6).Value = 1.036438012246326
$ws.Cells.Item(31,7).Value = 1.018633946610472
$ws.Cells.Item(31,8).Value = 1.036438012246326
$ws.Cells.Item(31,9).Value = 1.018633946610472
$ws.Cells.Item(31,10).Value = 0.9757885550309708
$ws.Cells.Item(31,11).Value = 1.036438012246326
$ws.Cells.Item(31,12).Value = 0.9757885550309708
$ws.Cells.Item(31,13).Value = 0.9972112508207215
$ws.Cells.Item(31,14).Value = 0.9972112508207215
$ws.Cells.Item(31,15).Value = 1.003009218938094
$ws.Cells.Item(31,16).Value = 1.01028683796259
$ws.Cells.Item(31,17).Value = 1.01028683796259
$ws.Cells.Item(31,18).Value = 1.016824631533523
$ws.Cells.Item(31,19).Value = 1.016824631533523
$ws.Cells.Item(31,20).Value = 0.9967856326866187

# ---------------------------------------------------------------------
# Copy the bold/bordered "index" column formatting (column A) down onto
# the two freshly-appended rows so they match the rest of the table.
# ---------------------------------------------------------------------
$ws.Range("A29").Copy() | Out-Null
$ws.Range("A30:A31").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0
